$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("breastcancer")

# Row 4 (bim, 0.05)
$ws.Range("C4").Value = 91.39784946236557
$ws.Range("D4").Value = 92.47311827956987
$ws.Range("F4").Value = 7.526881720430107
$ws.Range("G4").Value = 7.526881720430107

# Row 5 (bim, 0.2)
$ws.Range("F5").Value = 7.526881720430107
$ws.Range("G5").Value = 7.526881720430107

# Row 6 (bim, 0.4)
$ws.Range("F6").Value = 7.526881720430107
$ws.Range("G6").Value = 7.526881720430107

# Row 7 (bim, 1)
$ws.Range("F7").Value = 7.526881720430107
$ws.Range("G7").Value = 7.526881720430107

# Row 8 (boundary, 0.3)
$ws.Range("C8").Value = 67.74193548387096
$ws.Range("D8").Value = 91.39784946236557
$ws.Range("F8").Value = 7.526881720430107
$ws.Range("G8").Value = 7.526881720430107

# Row 9 (fgsm, 0.05)
$ws.Range("C9").Value = 91.39784946236557
$ws.Range("D9").Value = 92.47311827956987
$ws.Range("F9").Value = 7.526881720430107
$ws.Range("G9").Value = 7.526881720430107

# Row 10 (fgsm, 0.2)
$ws.Range("F10").Value = 7.526881720430107
$ws.Range("G10").Value = 7.526881720430107

# Row 11 (fgsm, 0.4)
$ws.Range("F11").Value = 7.526881720430107
$ws.Range("G11").Value = 7.526881720430107

# Row 12 (fgsm, 1)
$ws.Range("F12").Value = 7.526881720430107
$ws.Range("G12").Value = 7.526881720430107
